$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (America MG vs Sport Recife): update a couple of odds
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9

# Row 6 (Ind. Medellin vs Jaguares de Cordoba): update a couple of odds
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.7

# Row 7 (Cusco vs Grau): update one odd
$ws.Range("N7").Value = 9

# Row 8 (FC Cincinnati vs New York City) is removed entirely; the row
# that followed it (Seattle Sounders vs Houston Dynamo) shifts up to
# become the new row 8, shrinking the used range to A1:BD8.
$ws.Rows.Item(8).Delete()
